$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 81.166664
$ws.Range("I11").Value = 81.166664
$ws.Range("K11").Value = 81.166664
$ws.Range("M11").Value = 58.833336
$ws.Range("H92").Value = 1382.125
$ws.Range("I92").Value = 1378.1111
$ws.Range("K92").Value = 1378.1111
$ws.Range("M92").Value = -130.1111000000001
$ws.Range("H101").Value = 1992.5
$ws.Range("J101").Value = 1992.5
$ws.Range("L101").Value = 5977.5
$ws.Range("N101").Value = -9221.5
$ws.Range("H129").Value = 1277.8572
$ws.Range("I129").Value = 1157.5
$ws.Range("J129").Value = 2000
$ws.Range("K129").Value = 3472.5
$ws.Range("L129").Value = 6000
$ws.Range("M129").Value = 1527.5
$ws.Range("N129").Value = -16000
$ws.Range("H135").Value = 1601.5
$ws.Range("I135").Value = 1535
$ws.Range("K135").Value = 13815
$ws.Range("M135").Value = -11280
$ws.Range("H138").Value = 3288.182
$ws.Range("I138").Value = 873.2
$ws.Range("K138").Value = 2619.6
$ws.Range("M138").Value = 2520.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1999.5
$ws.Range("I45").Value = 1999.5
$ws.Range("K45").Value = 1999.5
$ws.Range("M45").Value = -1622.5
$ws.Range("H102").Value = 1750.2858
$ws.Range("I102").Value = 1840.1666
$ws.Range("J102").Value = 1211
$ws.Range("K102").Value = 1840.1666
$ws.Range("L102").Value = 1211
$ws.Range("M102").Value = -218.1666
$ws.Range("N102").Value = -4455
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H110").Value = 1097.5
$ws.Range("I110").Value = 1097.5
$ws.Range("K110").Value = 1097.5
$ws.Range("M110").Value = 947.5
$ws.Range("H122").Value = 7392.1113
$ws.Range("I122").Value = 7604.1333
$ws.Range("J122").Value = 6332
$ws.Range("K122").Value = 22812.3999
$ws.Range("L122").Value = 18996
$ws.Range("M122").Value = -20362.3999
$ws.Range("N122").Value = -23896

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 1000
$ws.Range("M36").Value = -466
$ws.Range("H97").Value = 8945.200000000001
$ws.Range("I97").Value = 8945.200000000001
$ws.Range("K97").Value = 8945.200000000001
$ws.Range("M97").Value = -7954.200000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 364.14285
$ws.Range("I7").Value = 216.33333
$ws.Range("J7").Value = 475
$ws.Range("K7").Value = 216.33333
$ws.Range("L7").Value = 475
$ws.Range("M7").Value = -103.33333
$ws.Range("N7").Value = -701
$ws.Range("H31").Value = 6085.9287
$ws.Range("I31").Value = 5118.6665
$ws.Range("K31").Value = 5118.6665
$ws.Range("M31").Value = -4823.6665
$ws.Range("H34").Value = 6085.9287
$ws.Range("I34").Value = 5118.6665
$ws.Range("K34").Value = 5118.6665
$ws.Range("M34").Value = -4916.6665
$ws.Range("H62").Value = 3626.2856
$ws.Range("I62").Value = 3476.8
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3476.8
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2852.8
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 3626.2856
$ws.Range("I65").Value = 3476.8
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 17384
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -14264
$ws.Range("N65").Value = -26240
$ws.Range("H107").Value = 767.7273
$ws.Range("I107").Value = 841.3333
$ws.Range("J107").Value = 679.4
$ws.Range("K107").Value = 841.3333
$ws.Range("L107").Value = 679.4
$ws.Range("M107").Value = 1078.6667
$ws.Range("N107").Value = -4519.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 3000
$ws.Range("N32").Value = -3566
$ws.Range("H35").Value = 7875
$ws.Range("J35").Value = 3833.3333
$ws.Range("L35").Value = 11499.9999
$ws.Range("N35").Value = -12075.9999
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H92").Value = 1190
$ws.Range("J92").Value = 275
$ws.Range("L92").Value = 825
$ws.Range("N92").Value = -3321
$ws.Range("H122").Value = 2263
$ws.Range("I122").Value = 1578.75
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 14208.75
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -11758.75
$ws.Range("N122").Value = -49900
$ws.Range("H131").Value = 1964.5
$ws.Range("I131").Value = 1619.6666
$ws.Range("J131").Value = 2999
$ws.Range("K131").Value = 4858.9998
$ws.Range("L131").Value = 8997
$ws.Range("M131").Value = 181.0002000000004
$ws.Range("N131").Value = -19077

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 662.5
$ws.Range("I13").Value = 662.5
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 662.5
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -523.5
$ws.Range("N13").ClearContents()
$ws.Range("H80").Value = 4000
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 5000
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 4000
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 25000
$ws.Range("N83").Value = -34984
$ws.Range("H93").Value = 40000
$ws.Range("J93").Value = 40000
$ws.Range("L93").Value = 40000
$ws.Range("N93").Value = -43744
$ws.Range("H122").Value = 3716.5
$ws.Range("I122").Value = 4059.8
$ws.Range("K122").Value = 12179.4
$ws.Range("M122").Value = -9729.400000000001
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 10002
$ws.Range("I21").Value = 10002
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 10002
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -9828
$ws.Range("N21").ClearContents()
$ws.Range("H40").Value = 6504.385
$ws.Range("I40").Value = 6414.273
$ws.Range("K40").Value = 6414.273
$ws.Range("M40").Value = -6278.273
$ws.Range("H93").Value = 3271.2856
$ws.Range("I93").Value = 3271.2856
$ws.Range("K93").Value = 3271.2856
$ws.Range("M93").Value = -2023.2856
$ws.Range("H132").Value = 8450
$ws.Range("I132").Value = 8450
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 25350
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -22820
$ws.Range("N132").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 999.5
$ws.Range("I10").Value = 999
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 999
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -830
$ws.Range("N10").Value = -1338
$ws.Range("H126").Value = 1571.5714
$ws.Range("I126").Value = 1571.5714
$ws.Range("K126").Value = 4714.7142
$ws.Range("M126").Value = -2244.7142
